$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 16.94169430348854
$ws.Range("C2").Value = 5.636622064413784
$ws.Range("D2").Value = 3.870082907719691
$ws.Range("F2").Value = 57.75520820565022
$ws.Range("G2").Value = 3.804516726942671
$ws.Range("I2").Value = 42.09939780854234
$ws.Range("J2").Value = 10.64191498200138
$ws.Range("K2").Value = 15.40375714176112
$ws.Range("L2").Value = 11.93954911668933
$ws.Range("M2").Value = 17.60402457809402

# Row 3
$ws.Range("B3").Value = 16.9114781145587
$ws.Range("C3").Value = 5.503319704475319
$ws.Range("D3").Value = 3.884552454888674
$ws.Range("F3").Value = 57.56375546146774
$ws.Range("G3").Value = 3.807580588487176
$ws.Range("I3").Value = 41.99529817583594
$ws.Range("J3").Value = 10.65082586175178
$ws.Range("K3").Value = 15.38256010984919
$ws.Range("L3").Value = 11.96143220058406
$ws.Range("M3").Value = 17.63456705616493

# Row 4
$ws.Range("B4").Value = 16.89780840594036
$ws.Range("C4").Value = 5.421737974753001
$ws.Range("D4").Value = 3.894117155200859
$ws.Range("F4").Value = 57.45274711121302
$ws.Range("G4").Value = 3.809560628227747
$ws.Range("I4").Value = 41.93525426547308
$ws.Range("J4").Value = 10.65669546977638
$ws.Range("K4").Value = 15.37360434875898
$ws.Range("L4").Value = 11.97636402433608
$ws.Range("M4").Value = 17.65640887744636

# Row 5
$ws.Range("B5").Value = 16.89347174508333
$ws.Range("C5").Value = 5.388615212229653
$ws.Range("D5").Value = 3.898186162793965
$ws.Range("F5").Value = 57.40917508727381
$ws.Range("G5").Value = 3.810392445375092
$ws.Range("I5").Value = 41.91176792023661
$ws.Range("J5").Value = 10.65918777481414
$ws.Range("K5").Value = 15.37097942068446
$ws.Range("L5").Value = 11.98282531822759
$ws.Range("M5").Value = 17.66608648770728

# Row 6
$ws.Range("B6").Value = 16.89282631733405
$ws.Range("C6").Value = 5.383124244145518
$ws.Range("D6").Value = 3.898872172410538
$ws.Range("F6").Value = 57.4020410619211
$ws.Range("G6").Value = 3.810532076532541
$ws.Range("I6").Value = 41.90792756138548
$ws.Range("J6").Value = 10.65960769061649
$ws.Range("K6").Value = 15.37060553348918
$ws.Range("L6").Value = 11.9839209592823
$ws.Range("M6").Value = 17.66774037431558

# Row 7
$ws.Range("B7").Value = 16.89774491746163
$ws.Range("C7").Value = 5.421290699555305
$ws.Range("D7").Value = 3.894171337286209
$ws.Range("F7").Value = 57.45215271909041
$ws.Range("G7").Value = 3.809571745329597
$ws.Range("I7").Value = 41.93493353430271
$ws.Range("J7").Value = 10.6567286750722
$ws.Range("K7").Value = 15.37356479499068
$ws.Range("L7").Value = 11.97644963884817
$ws.Range("M7").Value = 17.6565362472741

# Row 8
$ws.Range("B8").Value = 16.93026620106729
$ws.Range("C8").Value = 5.590636935005877
$ws.Range("D8").Value = 3.874930945534695
$ws.Range("F8").Value = 57.68784894018832
$ws.Range("G8").Value = 3.805552689342324
$ws.Range("I8").Value = 42.06270524799504
$ws.Range("J8").Value = 10.64490493111558
$ws.Range("K8").Value = 15.39560848408125
$ws.Range("L8").Value = 11.94678424904934
$ws.Range("M8").Value = 17.61391475446068

# Row 9
$ws.Range("B9").Value = 17.03246883009831
$ws.Range("C9").Value = 5.922484176218295
$ws.Range("D9").Value = 3.842588286238307
$ws.Range("F9").Value = 58.20107544803535
$ws.Range("G9").Value = 3.798451442505667
$ws.Range("I9").Value = 42.34362833928076
$ws.Range("J9").Value = 10.62486825035218
$ws.Range("K9").Value = 15.47083572156843
$ws.Range("L9").Value = 11.90045932362849
$ws.Range("M9").Value = 17.5548274979348

# Row 10
$ws.Range("B10").Value = 17.1304790945139
$ws.Range("C10").Value = 6.163144622472264
$ws.Range("D10").Value = 3.822096400751178
$ws.Range("F10").Value = 58.60796129860866
$ws.Range("G10").Value = 3.793704200709024
$ws.Range("I10").Value = 42.56800063695456
$ws.Range("J10").Value = 10.61205290765868
$ws.Range("K10").Value = 15.54527644017652
$ws.Range("L10").Value = 11.87362392728881
$ws.Range("M10").Value = 17.52632494786375

# Row 11
$ws.Range("B11").Value = 17.17991108943601
$ws.Range("C11").Value = 6.271361125786915
$ws.Range("D11").Value = 3.813481259319709
$ws.Range("F11").Value = 58.79923429210474
$ws.Range("G11").Value = 3.791645444465344
$ws.Range("I11").Value = 42.67384977281963
$ws.Range("J11").Value = 10.60663359750338
$ws.Range("K11").Value = 15.5832142304963
$ws.Range("L11").Value = 11.86297372784874
$ws.Range("M11").Value = 17.51658826839228

# Row 12
$ws.Range("B12").Value = 17.19931363999341
$ws.Range("C12").Value = 6.312114096507577
$ws.Range("D12").Value = 3.810320331158237
$ws.Range("F12").Value = 58.87252454455439
$ws.Range("G12").Value = 3.790880250259699
$ws.Range("I12").Value = 42.71446362296808
$ws.Range("J12").Value = 10.60464022930936
$ws.Range("K12").Value = 15.59815718132868
$ws.Range("L12").Value = 11.85916423943646
$ws.Range("M12").Value = 17.51336471656156

# Row 13
$ws.Range("B13").Value = 17.19510476664529
$ws.Range("C13").Value = 6.303347959701466
$ws.Range("D13").Value = 3.81099658501603
$ws.Range("F13").Value = 58.85670243671382
$ws.Range("G13").Value = 3.791044408826361
$ws.Range("I13").Value = 42.70569329187417
$ws.Range("J13").Value = 10.60506692478305
$ws.Range("K13").Value = 15.594913453669
$ws.Range("L13").Value = 11.85997474738393
$ws.Range("M13").Value = 17.51403836614869

# Row 14
$ws.Range("B14").Value = 17.18149372561765
$ws.Range("C14").Value = 6.274718677923993
$ws.Range("D14").Value = 3.813219175928531
$ws.Range("F14").Value = 58.80524687479885
$ws.Range("G14").Value = 3.791582203099725
$ws.Range("I14").Value = 42.67718052979189
$ws.Range("J14").Value = 10.60646842455475
$ws.Range("K14").Value = 15.58443209283308
$ws.Range("L14").Value = 11.86265584202381
$ws.Range("M14").Value = 17.51631378065106

# Row 15
$ws.Range("B15").Value = 17.17324520051989
$ws.Range("C15").Value = 6.257151670540878
$ws.Range("D15").Value = 3.814593781892063
$ws.Range("F15").Value = 58.77383987895528
$ws.Range("G15").Value = 3.791913492104727
$ws.Range("I15").Value = 42.65978441577737
$ws.Range("J15").Value = 10.60733453565229
$ws.Range("K15").Value = 15.57808677855268
$ws.Range("L15").Value = 11.86432718433585
$ws.Range("M15").Value = 17.51776787367351

# Row 16
$ws.Range("B16").Value = 17.12734501982037
$ws.Range("C16").Value = 6.156043129413349
$ws.Range("D16").Value = 3.822673625510077
$ws.Range("F16").Value = 58.59558302554402
$ws.Range("G16").Value = 3.793840766439294
$ws.Range("I16").Value = 42.5611581669204
$ws.Range("J16").Value = 10.61241531322106
$ws.Range("K16").Value = 15.5428783391972
$ws.Range("L16").Value = 11.87435124419538
$ws.Range("M16").Value = 17.52702617487294

# Row 17
$ws.Range("B17").Value = 17.10041903718423
$ws.Range("C17").Value = 6.093660336282246
$ws.Range("D17").Value = 3.82781122220591
$ws.Range("F17").Value = 58.48779038201355
$ws.Range("G17").Value = 3.795048843164959
$ws.Range("I17").Value = 42.50161380205066
$ws.Range("J17").Value = 10.61563717872913
$ws.Range("K17").Value = 15.52231639249479
$ws.Range("L17").Value = 11.88089925427826
$ws.Range("M17").Value = 17.53353238057231

# Row 18
$ws.Range("B18").Value = 17.08538897512197
$ws.Range("C18").Value = 6.057663186260315
$ws.Range("D18").Value = 3.830832762794576
$ws.Range("F18").Value = 58.4263742918637
$ws.Range("G18").Value = 3.795753188387876
$ws.Range("I18").Value = 42.46772214755145
$ws.Range("J18").Value = 10.61752895719696
$ws.Range("K18").Value = 15.51087380200751
$ws.Range("L18").Value = 11.8848121097107
$ws.Range("M18").Value = 17.53757859250858

# Row 19
$ws.Range("B19").Value = 17.08037893774884
$ws.Range("C19").Value = 6.045456630210957
$ws.Range("D19").Value = 3.831867236978232
$ws.Range("F19").Value = 58.40568092308839
$ws.Range("G19").Value = 3.795993300221192
$ws.Range("I19").Value = 42.45630865140156
$ws.Range("J19").Value = 10.61817612507559
$ws.Range("K19").Value = 15.50706577170084
$ws.Range("L19").Value = 11.88616212918319
$ws.Range("M19").Value = 17.53900080857343

# Row 20
$ws.Range("B20").Value = 17.10323814312681
$ws.Range("C20").Value = 6.100313454704425
$ws.Range("D20").Value = 3.827257431943808
$ws.Range("F20").Value = 58.49920492141345
$ws.Range("G20").Value = 3.794919259594796
$ws.Range("I20").Value = 42.50791557463745
$ws.Range("J20").Value = 10.61529020743068
$ws.Range("K20").Value = 15.52446555780155
$ws.Range("L20").Value = 11.88018703658661
$ws.Range("M20").Value = 17.53280832517018

# Row 21
$ws.Range("B21").Value = 17.18547317288056
$ws.Range("C21").Value = 6.283134274993597
$ws.Range("D21").Value = 3.812563595552047
$ws.Range("F21").Value = 58.82033753764831
$ws.Range("G21").Value = 3.791423849388339
$ws.Range("I21").Value = 42.68554110976211
$ws.Range("J21").Value = 10.60605517580818
$ws.Range("K21").Value = 15.58749514725716
$ws.Range("L21").Value = 11.86186227709276
$ws.Range("M21").Value = 17.51563286402112

# Row 22
$ws.Range("B22").Value = 17.24319691353707
$ws.Range("C22").Value = 6.4012808063016
$ws.Range("D22").Value = 3.803551503067357
$ws.Range("F22").Value = 59.03521062933452
$ws.Range("G22").Value = 3.789223364258753
$ws.Range("I22").Value = 42.80471888424599
$ws.Range("J22").Value = 10.60036222151222
$ws.Range("K22").Value = 15.63204544610552
$ws.Range("L22").Value = 11.85118856257048
$ws.Range("M22").Value = 17.50710894256157

# Row 23
$ws.Range("B23").Value = 17.2120292884982
$ws.Range("C23").Value = 6.338360165522388
$ws.Range("D23").Value = 3.808307394437773
$ws.Range("F23").Value = 58.92008179302833
$ws.Range("G23").Value = 3.790390148209427
$ws.Range("I23").Value = 42.74083314169972
$ws.Range("J23").Value = 10.60336937472178
$ws.Range("K23").Value = 15.60796419993375
$ws.Range("L23").Value = 11.85676628877143
$ws.Range("M23").Value = 17.51141147533532

# Row 24
$ws.Range("B24").Value = 17.10196222245937
$ws.Range("C24").Value = 6.09730598898327
$ws.Range("D24").Value = 3.827507589124549
$ws.Range("F24").Value = 58.49404267591607
$ws.Range("G24").Value = 3.794977813774871
$ws.Range("I24").Value = 42.50506547935655
$ws.Range("J24").Value = 10.6154469501467
$ws.Range("K24").Value = 15.52349273968872
$ws.Range("L24").Value = 11.88050856817383
$ws.Range("M24").Value = 17.53313471827916

# Row 25
$ws.Range("B25").Value = 17.0007535809014
$ws.Range("C25").Value = 5.833056032400104
$ws.Range("D25").Value = 3.850762523457111
$ws.Range("F25").Value = 58.05689827876006
$ws.Range("G25").Value = 3.800289575513133
$ws.Range("I25").Value = 42.26443265200996
$ws.Range("J25").Value = 10.62995301628857
$ws.Range("K25").Value = 15.44709018060732
$ws.Range("L25").Value = 11.91172533016048
$ws.Range("M25").Value = 17.5681918217838

Write-Host "Updated loading_percent values for 380 kV case"